# "Update countries & provincias Spain"
#
# This workbook ("Pais" sheet) is a COVID-19 daily country snapshot. The
# countries are listed sorted by total cases (column B) descending. This
# edit refreshes the case counters with a newer pull of the source data
# (timestamped ~30 minutes after the previous one) and re-sorts the few
# country rows whose case counts crossed a neighbour's count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp footer (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 23:16"

# --- Straight numeric refreshes (country keeps its row) -----------------
# China (row 4): no new cases/deaths today
$ws.Range("C4").Value = 0
$ws.Range("G4").Value = 0

# Espana (row 6)
$ws.Range("B6").Value = 25496
$ws.Range("C6").Value = 3925
$ws.Range("E6").Value = 21993

# Alemania (row 8)
$ws.Range("B8").Value = 22255
$ws.Range("C8").Value = 2407
$ws.Range("E8").Value = 21962

# Austria (row 15)
$ws.Range("B15").Value = 2992
$ws.Range("C15").Value = 343
$ws.Range("E15").Value = 2975

# Noruega (row 17)
$ws.Range("B17").Value = 2164
$ws.Range("C17").Value = 205
$ws.Range("E17").Value = 2151

# Jordania (row 76)
$ws.Range("B76").Value = 100
$ws.Range("C76").Value = 16
$ws.Range("E76").Value = 99

# Republica de Chipre (row 84)
$ws.Range("E84").Value = 80
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 1

# --- Canada overtakes Dinamarca & Portugal (rows 19-21) ------------------
# Canada's updated totals (1328) now exceed Dinamarca (1326) and Portugal
# (1280), so it moves from row 21 up to row 19; Dinamarca and Portugal
# each drop one place, keeping their own (unchanged) figures.
$ws.Range("A19").Value = "Canada"
$ws.Range("B19").Value = 1328
$ws.Range("C19").Value = 241
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 1295
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 7
$ws.Range("H19").Value = 19

$ws.Range("A20").Value = "Dinamarca"
$ws.Range("B20").Value = 1326
$ws.Range("C20").Value = 71
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1312
$ws.Range("F20").Value = 42
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 13

$ws.Range("A21").Value = "Portugal"
$ws.Range("B21").Value = 1280
$ws.Range("C21").Value = 260
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 1263
$ws.Range("F21").Value = 26
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 12

# --- Crucero & Luxemburgo overtake Pakistan (rows 30-32) ------------------
$ws.Range("A30").Value = "Crucero"
$ws.Range("B30").Value = 712
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 567
$ws.Range("E30").Value = 137
$ws.Range("F30").Value = 15
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 8

$ws.Range("A31").Value = "Luxemburgo"
$ws.Range("B31").Value = 670
$ws.Range("C31").Value = 186
$ws.Range("D31").Value = 6
$ws.Range("E31").Value = 656
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 8

$ws.Range("A32").Value = "Pakistan"
$ws.Range("B32").Value = 645
$ws.Range("C32").Value = 144
$ws.Range("D32").Value = 13
$ws.Range("E32").Value = 629
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 3

# --- Ecuador overtakes Grecia & Finlandia (rows 35-37) --------------------
$ws.Range("A35").Value = "Ecuador"
$ws.Range("B35").Value = 532
$ws.Range("C35").Value = 106
$ws.Range("D35").Value = 3
$ws.Range("E35").Value = 522
$ws.Range("F35").Value = 2
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 7

$ws.Range("A36").Value = "Grecia"
$ws.Range("B36").Value = 530
$ws.Range("C36").Value = 35
$ws.Range("D36").Value = 19
$ws.Range("E36").Value = 498
$ws.Range("F36").Value = 18
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 13

$ws.Range("A37").Value = "Finlandia"
$ws.Range("B37").Value = 523
$ws.Range("C37").Value = 73
$ws.Range("D37").Value = 10
$ws.Range("E37").Value = 512
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 1

# --- Ghana overtakes Bolivia (rows 117-118) -------------------------------
$ws.Range("A117").Value = "Ghana"
$ws.Range("B117").Value = 21
$ws.Range("C117").Value = 5
$ws.Range("E117").Value = 20
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 1

$ws.Range("A118").Value = "Bolivia"
$ws.Range("E118").Value = 19
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 0

# --- Polinesia Francesa overtakes Guam (rows 128-129) ---------------------
$ws.Range("A128").Value = "Polinesia Francesa"
$ws.Range("C128").Value = 4

$ws.Range("A129").Value = "Guam"
$ws.Range("C129").Value = 1

# --- Mayotte / Seychelles tie reorders (rows 137-138, figures unchanged) -
$ws.Range("A137").Value = "Mayotte"
$ws.Range("A138").Value = "Seychelles"
